# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Camote" (Vega Modelo de Temuco)
# at row 65, pushing the existing rows 65:146 down to 66:147.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 65:146 down one row (mirrors Excel's Rows.Insert UI action);
# this also extends the sheet dimension from R146 to R147 automatically.
$ws.Rows.Item(65).Insert()

# Populate the newly-opened row 65 with the new weekly record.
$ws.Range("A65").Value = 10
$ws.Range("B65").Value = "Vega Modelo de Temuco"
$ws.Range("C65").Value = "La Araucanía"
$ws.Range("D65").Value = 44893
$ws.Range("E65").Value = 9
$ws.Range("F65").Value = 100114002
$ws.Range("G65").Value = "Camote"
$ws.Range("H65").Value = "Sin especificar"
$ws.Range("I65").Value = "Primera"
$ws.Range("J65").Value = 40
$ws.Range("K65").Value = 22000
$ws.Range("L65").Value = 22000
$ws.Range("M65").Value = 22000
$ws.Range("N65").Value = "$/malla 20 kilos"
$ws.Range("O65").Value = "Perú"
$ws.Range("P65").Value = 1100
$ws.Range("Q65").Value = 20
$ws.Range("R65").Value = "Hortaliza"
